$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-49, replacing the old "Strike#" values
$gValues = @{
    2  = 1
    3  = 1
    4  = 0
    5  = 0
    6  = 0
    7  = 2
    8  = 1
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 3
    15 = 1
    16 = 2
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 3
    27 = 1
    28 = 1
    29 = 1
    30 = 1
    31 = 2
    32 = 0
    33 = 2
    34 = 1
    35 = 0
    36 = 2
    37 = 1
    38 = 2
    39 = 1
    40 = 2
    41 = 2
    42 = 1
    43 = 1
    44 = 2
    45 = 1
    46 = 0
    47 = 1
    48 = 3
    49 = 0
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
